$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Fix bug for swap scene: CanClone (col K) should be enabled for every
# scene row, and the Share flag (col J) for row 12 was incorrectly left on.
for ($r = 11; $r -le 35; $r++) {
    $ws.Cells.Item($r, 11).Value = 1
}
$ws.Cells.Item(12, 10).Value = 0
